$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 29: new submission entry (LogReg trained on all data, tile avg)
$ws.Range("I29").Value = "submission trained on all data"
$ws.Range("A29").Value = "2023-03-09-1757_LogReg_tile_avg.csv"

# Row 30: new submission entry (RF trained on all data, tile avg)
$ws.Range("A30").Value = "2023-03-09-1800_RF_tile_avg.csv"

$ws.Range("N29").Value = "March 9"
$ws.Range("O29").Value = 0.587

$ws.Range("I30").Value = "submission trained on all data"
$ws.Range("N30").Value = "March 9"
$ws.Range("O30").Value = 0.633

# Update the view: scroll so column F is the leftmost visible column,
# and select N31 as the active cell.
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("N31").Select()
